$d = $word.ActiveDocument

function Touch-Font([object]$rng) {
    # Forces Word to split this sub-range into its own run by re-asserting
    # character formatting that already matches its effective font, mirroring
    # how Word itself fragments runs around proofing marks / edits without
    # altering visible formatting or text.
    $rng.Font.NameAscii = "Times New Roman"
    $rng.Font.NameFarEast = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Name = "Times New Roman"
}

# ---------------------------------------------------------------------------
# 1) " A 82" -> " " / "A" / " 82"  (Group ID line)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("A 82", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $a = $d.Range($rng.Start, $rng.Start + 1)
    Touch-Font $a
}

Write-Output "step1 done"

# ---------------------------------------------------------------------------
# 2) "Aarthi Kalyanapu " -> "Aarthi " / "Kalyanapu" / " "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Aarthi Kalyanapu ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $mid = $d.Range($rng.Start + 7, $rng.Start + 16)
    Touch-Font $mid
}

Write-Output "step2 done"

# ---------------------------------------------------------------------------
# 3) "Kingsley Anyaeche " -> "Kingsley " / "Anyaeche" / " "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Kingsley Anyaeche ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $mid = $d.Range($rng.Start + 9, $rng.Start + 17)
    Touch-Font $mid
}

Write-Output "step3 done"

# ---------------------------------------------------------------------------
# 5) Remove the first (highlighted, ListParagraph) of the two empty
#    paragraphs directly after the "1. Introduction" heading.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "1. Introduction" + [char]13) {
        $target = $paras.Item($i + 1)
        if ($target.Style.NameLocal -eq "List Paragraph" -and $target.Range.Text -eq [char]13) {
            $target.Range.Delete()
        }
        break
    }
}

Write-Output "step5 done"

# ---------------------------------------------------------------------------
# 6) "...By comparing 40 balls with and without a coating, this study..."
#    -> "...By comparing 40 balls with coating and without coating, this study..."
# ---------------------------------------------------------------------------
$old6 = "By comparing 40 balls with and without a coating, this study aims to find out"
$new6 = "By comparing 40 balls with coating and without coating, this study aims to find out"
$null = $d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2)

Write-Output "step6 done"

# ---------------------------------------------------------------------------
# 7) Rewrite the "data set contains measurements" paragraph.
# ---------------------------------------------------------------------------
$old7 = 'The data set contains measurements of how far golf balls travelled when hit under similar conditions. It includes two groups of 40 balls each: Current balls, which have no special coating, and New balls, which have a coating applied. Each value represents the driving distance of one shot. By comparing the distances from the coated and uncoated balls, the dataset allows us to see whether the coating might help the ball travel farther'
$new7 = 'The data set contains measurements of how far golf balls travelled when hit under similar conditions. It includes two groups of 40 balls each: Current balls, which are without coating, and New balls, which are with coating applied. Each value represents the driving distance of one shot. By comparing the distances from the with coating and without coating balls, the dataset allows us to see whether the coating might help the ball travel farther'
$null = $d.Content.Find.Execute($old7, $false, $false, $false, $false, $false, $true, 1, $false, $new7, 2)

Write-Output "step7 done"

# ---------------------------------------------------------------------------
# 8) Rewrite the "This study aims to determine" paragraph.
# ---------------------------------------------------------------------------
$old8 = 'This study aims to determine whether adding a coating to golf balls affects how far they travel when hit. By comparing the average driving distances of 40 coated and 40 uncoated balls, the research question asks: Is there a significant difference in the mean driving distance between golf balls with a coating and those without?'
$new8 = 'This study aims to determine whether adding a coating to golf balls affects how far they travel when hit. By comparing the average driving distances of 40 with coating and 40 without coating balls, the research question asks: Is there a significant difference in the mean driving distance between golf balls with a coating and those without?'
$null = $d.Content.Find.Execute($old8, $false, $false, $false, $false, $false, $true, 1, $false, $new8, 2)

Write-Output "step8 done"
